# 2.6 - Node engine and OpenID upgrade
#
# The real content change behind this revision is the addition of two
# bookmarks around the document title / diagram heading:
#   - "_Hlk149316668" : a Word "Smart Lookup"-style bookmark that wraps the
#                        title paragraph through the "Ροή Δεδομένων" heading
#                        paragraph (i.e. from the very start of the document
#                        up to - and including - the end of that heading).
#   - "_GoBack"        : Word's standard "last edit location" bookmark,
#                         placed (collapsed / zero-length) immediately after
#                         the text "Ευκαρδία " in the title paragraph.

$d = $word.ActiveDocument

# --- Locate the two anchor paragraphs -------------------------------------
# Paragraph 1: the title ("Ευκαρδία" / "Προσχέδιο ροής καταχώρισης στοιχείων")
# Paragraph 3: the heading paragraph that holds the diagram + "Ροή Δεδομένων"
$titlePara = $d.Paragraphs.Item(1)
$headingPara = $d.Paragraphs.Item(3)

# --- Bookmark 1: "_Hlk149316668" -------------------------------------------
# Spans from the very beginning of the document through the end of the
# "Ροή Δεδομένων" heading paragraph.
$hlkRange = $d.Range($titlePara.Range.Start, $headingPara.Range.End)
$d.Bookmarks.Add("_Hlk149316668", $hlkRange)

# --- Bookmark 2: "_GoBack" ---------------------------------------------------
# Collapsed bookmark right after "Ευκαρδία " (before the line break) marking
# the last edit position.
$afterTitleWord = $d.Range(0, 0)
$afterTitleWord.Find.Execute("Ευκαρδία ") | Out-Null
$goBackPos = $afterTitleWord.End
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
